$d = $word.ActiveDocument

# The edit renames the cosmetic "name" label stored on each picture's
# wp:docPr / pic:cNvPr elements (footer "PearsonLogo.png" pictures:
# image1.png -> image2.png; header "BTec_Logo-Orange" picture:
# image2.jpg -> image1.jpg). The embedded relationship (r:embed) and the
# actual media parts are left untouched - only the display "name"
# attribute text changes, in both places it is duplicated (wp:docPr and
# pic:cNvPr) for every affected picture.
#
# InlineShape.Name only ever reaches the wp:docPr copy of the attribute
# in this object model, so the rename is done via a full WordOpenXML
# round-trip: pull the whole package as flat OOXML text, patch the
# "name=" attribute values with ordinary string replacement, and write
# the patched text back. A plain get/set round trip of
# Content.WordOpenXML is a byte-for-byte no-op otherwise, so this only
# changes the targeted attribute values.

$xml = $d.Content.WordOpenXML

$xml = $xml.Replace('name="image1.png"', 'name="image2.png"')
$xml = $xml.Replace('name="image2.jpg"', 'name="image1.jpg"')

$d.Content.WordOpenXML = $xml
